# Update the cryptos worksheet with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "0.380", "89.313.33")
# but must stay plain text so formatting like trailing zeros / multi-dot
# separators survive - force text format on each Price cell right before
# writing its refreshed value (the engine doesn't support comma-joined
# multi-area Range() for NumberFormat, so this is done cell-by-cell).
function Set-PriceText([string]$cellRef, [string]$value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-PriceText "D2" "89.313.33"
$ws.Range("E2").Value = "  +2.78%  "

Set-PriceText "D3" "3.271.30"
$ws.Range("E3").Value = "  -1.45%  "

$ws.Range("E4").Value = "  -0.07%  "

Set-PriceText "D5" "212.54"
$ws.Range("E5").Value = "  -2.80%  "

Set-PriceText "D6" "625.93"
$ws.Range("E6").Value = "  -1.89%  "

Set-PriceText "D7" "0.380"
$ws.Range("E7").Value = "  +17.73%  "

Set-PriceText "D8" "0.725"
$ws.Range("E8").Value = "  +18.35%  "

Set-PriceText "D9" "0.999"
$ws.Range("E9").Value = "  +0.01%  "

Set-PriceText "D10" "3.266.49"
$ws.Range("E10").Value = "  -1.64%  "

$ws.Range("E11").Value = "  -4.26%  "

Set-PriceText "D12" "0.186"
$ws.Range("E12").Value = "  +11.20%  "

Set-PriceText "D13" "0.0000263"
$ws.Range("E13").Value = "  -4.30%  "

Set-PriceText "D14" "34.09"
$ws.Range("E14").Value = "  -0.53%  "

Set-PriceText "D15" "3.868.93"
$ws.Range("E15").Value = "  -1.59%  "

$ws.Range("E16").Value = "  +0.34%  "

Set-PriceText "D17" "88.926.89"
$ws.Range("E17").Value = "  +2.46%  "

Set-PriceText "D18" "3.286.82"
$ws.Range("E18").Value = "  -1.10%  "

Set-PriceText "D19" "14.07"
$ws.Range("E19").Value = "  -4.27%  "

Set-PriceText "D20" "3.08"
$ws.Range("E20").Value = "  -3.72%  "

Set-PriceText "D21" "436.79"
$ws.Range("E21").Value = "  -2.36%  "

Set-PriceText "D22" "8.89"
$ws.Range("E22").Value = "  -2.89%  "

Set-PriceText "D23" "5.34"
$ws.Range("E23").Value = "  +1.66%  "

Set-PriceText "D24" "7.41"
$ws.Range("E24").Value = "  -0.57%  "

Set-PriceText "D25" "5.26"
$ws.Range("E25").Value = "  -2.41%  "

Set-PriceText "D26" "12.18"
$ws.Range("E26").Value = "  -0.45%  "

Set-PriceText "D27" "3.464.59"
$ws.Range("E27").Value = "  +0.83%  "

Set-PriceText "D28" "76.82"
$ws.Range("E28").Value = "  -2.35%  "

$ws.Range("E29").Value = "  +3.20%  "

$ws.Range("E30").Value = "  +0.03%  "

Set-PriceText "D31" "0.184"
$ws.Range("E31").Value = "  +6.06%  "

Set-PriceText "D32" "0.999"
$ws.Range("E32").Value = "  +0.06%  "

Set-PriceText "D33" "8.84"
$ws.Range("E33").Value = "  -4.66%  "

Set-PriceText "D34" "561.20"
$ws.Range("E34").Value = "  -7.48%  "

$ws.Range("E35").Value = "  -11.65%  "

# --- Row 36/37 swap: RenderToken <-> PancakeSwap ---
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-PriceText "D36" "1.96"
$ws.Range("E36").Value = "  -4.41%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-PriceText "D37" "7.11"
$ws.Range("E37").Value = "  +9.23%  "

Set-PriceText "D38" "0.139"
$ws.Range("E38").Value = "  -8.25%  "

Set-PriceText "D39" "22.65"
$ws.Range("E39").Value = "  -2.81%  "

Set-PriceText "D40" "21.82"
$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("E41").Value = "  +0.01%  "

Set-PriceText "D42" "3.10"
$ws.Range("E42").Value = "  -1.72%  "

Set-PriceText "D43" "2.03"
$ws.Range("E43").Value = "  -1.44%  "

Set-PriceText "D44" "0.398"
$ws.Range("E44").Value = "  -4.85%  "

$ws.Range("E45").Value = "  +0.11%  "

Set-PriceText "D46" "155.17"
$ws.Range("E46").Value = "  -0.91%  "

Set-PriceText "D47" "180.57"
$ws.Range("E47").Value = "  -4.54%  "

Set-PriceText "D48" "44.95"
$ws.Range("E48").Value = "  -1.50%  "

# --- Row 49/50 swap: ImmutableX <-> Stellar ---
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-PriceText "D49" "0.131"
$ws.Range("E49").Value = "  +16.57%  "

$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-PriceText "D50" "1.31"
$ws.Range("E50").Value = "  -4.38%  "

Set-PriceText "D51" "4.23"
$ws.Range("E51").Value = "  -0.86%  "
